$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.161638975143433
$ws.Range("B1").Value = 2.15468692779541
$ws.Range("C1").Value = 10.32306385040283
$ws.Range("D1").Value = 2.552456378936768
$ws.Range("E1").Value = 1.263031125068665
